$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-09-25 12:41:25"

# ------------------------------------------------------------------
# 1. Insert two brand-new rows right after the existing row 2. This
#    pushes the former rows 3-9 down to rows 5-11 (Excel's native
#    row-insert semantics - existing hyperlink ref/rId pairs for F2..F9
#    stay put, which is exactly what the target workbook shows).
# ------------------------------------------------------------------
$ws.Range("A3:A4").EntireRow.Insert()

# ------------------------------------------------------------------
# 2. New row 3
# ------------------------------------------------------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【低予算希望】LINE公式アカウント+社食注文システム開発依頼(社内利用のみ)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5400375"
$ws.Range("G3").Value = 118
$ws.Range("H3").Value = "◆開発,システム開発"

# ------------------------------------------------------------------
# 3. New row 4
# ------------------------------------------------------------------
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【急募】住宅展示場マッチング診断サービスのMVP開発依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5399759"
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = "◆開発"

# ------------------------------------------------------------------
# 4. Refresh the timestamp on all the pre-existing rows (now rows
#    2, 5-11) so the whole sheet reflects the same scrape run.
# ------------------------------------------------------------------
$ws.Range("A2").Value = $timestamp
$ws.Range("A5").Value = $timestamp
$ws.Range("A6").Value = $timestamp
$ws.Range("A7").Value = $timestamp
$ws.Range("A8").Value = $timestamp
$ws.Range("A9").Value = $timestamp
$ws.Range("A10").Value = $timestamp
$ws.Range("A11").Value = $timestamp

# ------------------------------------------------------------------
# 5. Brand-new row 12, appended at the bottom of the table.
# ------------------------------------------------------------------
$ws.Range("A12").Value = $timestamp
$ws.Range("B12").Value = "【SalesIQ活用】CRMと連携したリード獲得方法を教えてください"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "~ 5,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5400402"
$ws.Range("G12").Value = 10

# ------------------------------------------------------------------
# 6. F10/F11 kept their old (now stale) hyperlink relationships
#    because the insert above only shifted cell content, not the
#    hyperlink ref/rId pairs. Re-create correct hyperlinks for F10,
#    F11 (content unchanged from before, but now mapped to new
#    relationship ids) and the brand new F12, in that order, so
#    the new relationship ids come out as rId9, rId10, rId11.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5400231")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5399347")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5400402")

# Hyperlinks.Add() registers its own style xf; snap F10:F12 back onto
# the shared "Hyperlink" cell style used everywhere else in column F.
$ws.Range("F10:F12").Style = "Hyperlink"
